# "Update Work Week and Social Spending"
#
# Refreshes the Latvia GDP-per-Capita series on the "Data" sheet:
#   - the 1973 figure is revised (7846 -> 12506)
#   - the 1980-2010 figures are revised with newer source data
#   - six new rows are appended for 2011-2016
# The 1974-1979 rows are left untouched (still blank).
#
# The "Data" column has always stored its figures as text (even though
# they look numeric), so every cell written below is first forced to
# Text format - otherwise Excel would silently re-interpret a literal
# like "12506" as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Force Text format on every "Data" (column E) cell that is about to be
# (re)written, so the numeric-looking strings stay strings.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E9:E45").NumberFormat = "@"

# Revised figures for existing rows (year -> new Data value)
$revisedValues = [ordered]@{
    1973 = "12506"
    1980 = "13739"
    1981 = "13970"
    1982 = "14190"
    1983 = "14376"
    1984 = "14682"
    1985 = "14273"
    1986 = "14539"
    1987 = "14360"
    1988 = "14875"
    1989 = "15661"
    1990 = "15806"
    1991 = "13774.9340381701"
    1992 = "9411.19496932769"
    1993 = "8439.33034911096"
    1994 = "8704.26615324039"
    1995 = "8709.07005121388"
    1996 = "8947.26813728627"
    1997 = "9773.4995025883"
    1998 = "10429.7666358074"
    1999 = "10707.7057375105"
    2000 = "11309.7403901285"
    2001 = "12103.8281717161"
    2002 = "13015.4745034051"
    2003 = "14142.1245473969"
    2004 = "15371.1813266083"
    2005 = "17069.6571716176"
    2006 = "19128.6892442587"
    2007 = "21042.3664707057"
    2008 = "20342.1982902556"
    2009 = "17582.1116698334"
    2010 = "17140.226514887"
}

# Row numbers on the "Data" sheet line up with year - 1971
foreach ($year in $revisedValues.Keys) {
    $row = $year - 1971
    $ws.Range("E$row").Value = $revisedValues[$year]
}

# New rows for years 2011-2016
$newYears = 2011, 2012, 2013, 2014, 2015, 2016
$newValues = @{
    2011 = "18428"
    2012 = "19405"
    2013 = "20129"
    2014 = "20751"
    2015 = "21488"
    2016 = "22092"
}

foreach ($year in $newYears) {
    $row = $year - 1971
    $ws.Range("A$row").Value = 428
    $ws.Range("B$row").Value = "Latvia"
    $ws.Range("C$row").Value = "GDP per Capita"
    $ws.Range("D$row").Value = $year
    $ws.Range("E$row").Value = $newValues[$year]
}
